# Applies the "Modif url canonique termino" change:
#  - Update the Date metadata value on the "Metadata" sheet
#  - Update three canonical terminology URLs on the "Elements" sheet
#  - Widen column Z on "Elements" to fit the new (longer) URL text

$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsElements = $wb.Worksheets.Item("Elements")

# 1. Update the Date value (row 8, column B on Metadata sheet)
$wsMetadata.Range("B8").Value = "2025-07-25T07:22:51+00:00"

# 2. Update the canonical terminology URLs on the Elements sheet (column Z)
$wsElements.Range("Z3").Value = "https://mos.esante.gouv.fr/NOS/TRE_R14-TypeDiplome/FHIR/TRE-R14-TypeDiplome?vs"
$wsElements.Range("Z4").Value = "https://mos.esante.gouv.fr/NOS/TRE_R16-LieuFormation/FHIR/TRE-R16-LieuFormation?vs"
$wsElements.Range("Z7").Value = "https://mos.esante.gouv.fr/NOS/TRE_R50-DESCGroupe1Diplome/FHIR/TRE-R50-DESCGroupe1Diplome?vs"

# 3. Resize column Z to best fit the new (longer) content, matching the
#    canonical width recorded in the target workbook (~83.52 characters).
$wsElements.Columns.Item(26).ColumnWidth = 82.69
